$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old 2000-2009 rows (original rows 2-11); this shifts the
# 2010-2020 data that was in rows 12-22 up to rows 2-12.
$ws.Rows("2:11").Delete()

# Append the two new years that were not present before.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 1040
$ws.Range("C13").Value = 3.96

$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Value = 1203
# 2022 has no registered-unemployment-rate figure yet, so leave it blank.
$ws.Range("C14").ClearContents()

# Match the formatting used by the rest of column A (bold, centered, boxed)
# by copying the existing style rather than re-building it from scratch.
$ws.Range("A12").Copy()
$ws.Range("A13:A14").PasteSpecial(-4122)
